# Update description text for PSC statistics: change "mili(seconds|secods)" to "seconds"/"secods"
# Rows (1-indexed, header in row 1):
#   Row 6 (psc_start_ms)    -> column C: "the vector denoting the starting time of PSC in seconds"
#   Row 8 (psc_risetime_ms) -> column C: "rise time of the PSC in secods"
#   Row 9 (psc_decay_ms)    -> column C: "decay time of the PSC in seconds"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = "rise time of the PSC in secods"
$ws.Range("C9").Value = "decay time of the PSC in seconds"
$ws.Range("C6").Value = "the vector denoting the starting time of PSC in seconds"

# Scroll/selection state change observed in the diff
$ws.Range("J11").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
